$d = $word.ActiveDocument

# Use an existing run elsewhere in the document that already carries the
# exact character formatting we need (Times New Roman / sz 24 / szCs 24 /
# lang en-US) so the new run's <w:rPr> is produced faithfully, including
# the w:cs font and w:lang that plain Font.* property assignments don't
# reliably serialize.
$srcPara = $d.Paragraphs.Item(3)
$srcRange = $srcPara.Range
# Exclude the trailing paragraph mark from the source range so only the
# run-level (character) formatting is captured, not paragraph formatting.
$srcText = $d.Range($srcRange.Start, $srcRange.End - 1)

# Target: the last paragraph in the document, which currently only holds
# the (hidden) "_GoBack" bookmark.
$lastPara = $d.Paragraphs.Last
$targetRange = $lastPara.Range
$insertPos = $targetRange.Start

# Insert a same-length placeholder *before* the bookmark so the new run
# ends up ahead of <w:bookmarkStart>/<w:bookmarkEnd> in the XML, matching
# the diff ordering.
$placeholder = $srcText.Text
$targetRange.InsertBefore($placeholder)

$newRange = $d.Range($insertPos, $insertPos + $placeholder.Length)

# Clone the formatting from the source run onto the freshly inserted text.
$newRange.FormattedText = $srcText.FormattedText

# Finally, swap in the real note text (keeps the run's formatting).
$newRange.Text = "Notes: This file must be updated as soon as possible."
